# "Adding new CSV utility"
#
# Replace the old two-column (id / name / score) sample data on Sheet1
# with a small CSV-style contact table: name, email, country, state -
# one row for a US contact and one row for an Indian contact - and
# right-align the populated cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - US contact
$ws.Range("A1").Value = "John Wick"
$ws.Range("B1").Value = "john@test.com"
$ws.Range("C1").Value = "US"
$ws.Range("D1").Value = "California"

# Row 2 - Indian contact
$ws.Range("A2").Value = "Sachin Taware"
$ws.Range("B2").Value = "Sachin@Test.com"
$ws.Range("C2").Value = "IND"
$ws.Range("D2").Value = "Maharashtra"

# Right-align every populated cell (A1:D2), like the new CSV import style.
$ws.Range("A1:D2").HorizontalAlignment = -4152  ## xlRight
